# Updated cryptos list - apply scraped Price (D) / Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.522.24'
$ws.Range("E2").Value = '  +4.05%  '
$ws.Range("D3").Value = '1.737.24'
$ws.Range("E3").Value = '  +4.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.67'
$ws.Range("E5").Value = '  +3.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4802'
$ws.Range("E7").Value = '  +3.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2669'
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06232'
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("D10").Value = '1.737.24'
$ws.Range("E10").Value = '  +4.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07135'
$ws.Range("E11").Value = '  +2.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.76'
$ws.Range("E12").Value = '  +7.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6177'
$ws.Range("E13").Value = '  +7.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.535'
$ws.Range("E14").Value = '  +4.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.89'
$ws.Range("E15").Value = '  +2.49%  '
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = '26.537.61'
$ws.Range("E17").Value = '  +4.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006903'
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.75'
$ws.Range("E20").Value = '  +3.30%  '
$ws.Range("D21").Value = '1.960.60'
$ws.Range("E21").Value = '  +4.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.564'
$ws.Range("E22").Value = '  +3.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.891'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.343'
$ws.Range("E24").Value = '  +2.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.72'
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.36'
$ws.Range("E26").Value = '  +3.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.801'
$ws.Range("E27").Value = '  +5.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.411'
$ws.Range("E28").Value = '  +3.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.63'
$ws.Range("E29").Value = '  +2.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.992'
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.715'
$ws.Range("E31").Value = '  +3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07890'
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04580'
$ws.Range("E33").Value = '  +5.57%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9966'
$ws.Range("E35").Value = '  +5.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6343'
$ws.Range("E36").Value = '  +5.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9300'
$ws.Range("E37").Value = '  +2.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '111.42'
$ws.Range("E38").Value = '  +3.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.439'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.981'
$ws.Range("E40").Value = '  +8.38%  '
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01510'
$ws.Range("E42").Value = '  +3.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.713'
$ws.Range("E43").Value = '  +14.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3900'
$ws.Range("E44").Value = '  +5.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.909'
$ws.Range("E45").Value = '  +12.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1194'
$ws.Range("E46").Value = '  +8.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05334'
$ws.Range("E47").Value = '  +1.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.885'
$ws.Range("E48").Value = '  +4.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.83'
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.254'
$ws.Range("E50").Value = '  +5.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3440'
$ws.Range("E51").Value = '  +3.98%  '
